# M365 LogSize Estimator: switch license flags from yes/no dropdown to numeric counts
$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false
$ws = $wb.Worksheets("Sheet1")

# --- Fill in the measured Defender usage numbers (Count / MB columns) ---
$ws.Range("B8").Value  = 541458498
$ws.Range("C8").Value  = 486359.75
$ws.Range("B9").Value  = 12496671
$ws.Range("C9").Value  = 4397.0200000000004
$ws.Range("B10").Value = 11534952
$ws.Range("C10").Value = 23771.32
$ws.Range("B11").Value = 16119836
$ws.Range("C11").Value = 9387.42
$ws.Range("B12").Value = 25689
$ws.Range("C12").Value = 41.9

# --- Insert a new explanatory row under the licensing header ---
$ws.Rows(16).Insert()
$ws.Range("A16").Value = "(only licensed users count)"

# --- Replace the yes/no dropdown answers with numeric license counts ---
$ws.Range("B17").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("B19").Value = 0
$ws.Range("B20").Value = 0
$ws.Range("B21").Value = 0
$ws.Range("B22").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0

# --- Update the license-benefit formula to sum the counts instead of nested IFs ---
$ws.Range("B27").Formula = "=B17*5+B18*5+B19*5+B20*5+B21*5+B22*5+B23*5+B24*5"

# --- Replace the list-based data validation (pointing at Sheet2) with a numeric one ---
$ws.Range("B17:B24").Validation.Delete()
$ws.Range("B17:B24").Validation.Add(1, 1, 1, 0, 1000000)

# --- Update the instructional text above the license section ---
$ws.Range("A15").Value = "If you hold one of the following licenses, please insert the number of licenses for each"

# --- Remove the now unused helper Sheet2 (and its Table1/dropdown list) ---
$wb.Worksheets("Sheet2").Delete()

# --- Tidy up row heights to match the refreshed layout ---
$ws.Rows(8).RowHeight = 15.5
$ws.Rows(9).RowHeight = 15.5
$ws.Rows(10).RowHeight = 15.5
$ws.Rows(11).RowHeight = 15.5
$ws.Rows(12).RowHeight = 15.5
$ws.Rows(18).RowHeight = 15.5
$ws.Rows(19).RowHeight = 15.5
$ws.Rows(20).RowHeight = 15.5
$ws.Rows(21).RowHeight = 15.5
$ws.Rows(22).RowHeight = 15.5
$ws.Rows(23).RowHeight = 15.5
$ws.Rows(24).RowHeight = 15.5
$ws.Rows(27).RowHeight = 15.5
$ws.Rows(28).RowHeight = 15.5
$ws.Rows(29).RowHeight = 15.5

# --- Update the view to where the author left it ---
$ws.Range("D21").Select()
$excel.ActiveWindow.ScrollRow = 8
